# Update database and change read_price algorithm
# -----------------------------------------------------------------
# The workbook holds quarterly income-statement data for one ticker.
# Five new (older) quarters are inserted to the left of the existing
# quarterly columns (D:H -> I:M), new quarter data + headers fill the
# freed-up D:H columns, and the "EPS based on latest capital" row is
# recomputed for every period (including the previously existing ones).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$xlShiftToRight = [Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftToRight

# ------------------------------------------------------------------
# 1) Make room: insert 5 blank columns before column D. Everything
#    that used to live in D:H (values, styles, formats) slides right
#    into I:M automatically.
# ------------------------------------------------------------------
$ws.Range("D:H").Insert($xlShiftToRight)

# ------------------------------------------------------------------
# 2) Column widths for the freshly-inserted D:H block. The host's
#    ColumnWidth property reads/writes ~5/6 of a character narrower
#    than the value actually stored in the OOXML <col width=.../>
#    attribute (confirmed against this sheet's untouched I:M columns,
#    which still carry their original width="29"/"31" after the
#    shift), so compensate by that fixed offset.
# ------------------------------------------------------------------
$widthOffset = 5 / 6
$ws.Columns.Item(4).ColumnWidth = 29 - $widthOffset   # D
$ws.Columns.Item(5).ColumnWidth = 29 - $widthOffset   # E
$ws.Columns.Item(6).ColumnWidth = 31 - $widthOffset   # F
$ws.Columns.Item(7).ColumnWidth = 29 - $widthOffset   # G
$ws.Columns.Item(8).ColumnWidth = 29 - $widthOffset   # H

# ------------------------------------------------------------------
# 3) Header row 8 - period labels for the new quarters.
# ------------------------------------------------------------------
$ws.Range("D8").Value = "فصل دوم منتهی به 1399/06"
$ws.Range("E8").Value = "فصل سوم منتهی به 1399/09"
$ws.Range("F8").Value = "فصل چهارم منتهی به 1399/12"
$ws.Range("G8").Value = "فصل اول منتهی به 1400/03"
$ws.Range("H8").Value = "فصل دوم منتهی به 1400/06"

# ------------------------------------------------------------------
# 4) Header row 9 - publish dates for the new quarters.
# ------------------------------------------------------------------
$ws.Range("D9").Value = "1400-09-14 (5)"
$ws.Range("E9").Value = "1400-10-29 (2)"
$ws.Range("F9").Value = "1401-03-15 (9)"
$ws.Range("G9").Value = "1401-04-29 (2)"
$ws.Range("H9").Value = "1401-09-13 (5)"

# ------------------------------------------------------------------
# 5) Financial data rows 11-27 for the new D:H quarters.
# ------------------------------------------------------------------
$data = @{
    11 = @(26138536, 33894702, 57776495, 52463805, 50134459)
    12 = @(-13874333, -18607437, -33413838, -32047835, -19811168)
    13 = @(12264203, 15287265, 24362657, 20415970, 30323291)
    14 = @(-485432, -635080, -833378, -954441, -1063848)
    15 = @(0, 0, 0, 0, 0)
    16 = @(7148068, 113420, -2377223, 135721, -1857112)
    17 = @(18926839, 14765605, 21152056, 19597250, 27402331)
    18 = @(-458328, -697548, -867037, -771051, -135382)
    19 = @(137271, 457222, -1126329, -54433, 11235)
    20 = @(18605782, 14525279, 19158690, 18771766, 27278184)
    21 = @(-1593623, -1960628, -884930, -2162963, -2273947)
    22 = @(17012159, 12564651, 18273760, 16608803, 25004237)
    23 = @(0, 0, 0, 0, 0)
    24 = @(17012159, 12564651, 18273760, 16608803, 25004237)
    25 = @(680, 503, 731, 664, 1000)
    26 = @(25000000, 25000000, 25000000, 25000000, 25000000)
    27 = @(347, 256, 373, 339, 510)
}

$cols = @("D", "E", "F", "G", "H")
foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt 5; $i++) {
        $ws.Range($cols[$i] + $row).Value = $vals[$i]
    }
}

# ------------------------------------------------------------------
# 6) The "EPS based on latest capital" row (27) was recomputed with a
#    new read_price algorithm, so the shifted-right I:M values (which
#    the Insert step copied verbatim from the old D:H) must be
#    overwritten with the recalculated figures.
# ------------------------------------------------------------------
$ws.Range("I27").Value = 430
$ws.Range("J27").Value = 358
$ws.Range("K27").Value = 247
$ws.Range("L27").Value = 265
$ws.Range("M27").Value = 123
